$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "CSS" column (F) selectors, and the one missing
# Xpath (D8) value, for the homework table. Cells are written in the same
# order the original author typed them in (new shared strings are appended
# on first use), so that F4's value - added last by the author - lands at
# the end of the shared-strings table, matching row order elsewhere.

$ws.Range("F2").Value  = 'a[href*="+38"]'
$ws.Range("F3").Value  = 'nav a[href*="contacts"]'
$ws.Range("F5").Value  = "#fat-menu"
$ws.Range("F6").Value  = "input[name = 'search']"
$ws.Range("F7").Value  = ".button.button_color_green.button_size_medium.search-form__submit.ng-star-inserted"
$ws.Range("D8").Value  = "//button[@class = 'header__button']"
$ws.Range("F8").Value  = ".city-toggle__text"
$ws.Range("F9").Value  = ".top-information__inner.ng-star-inserted"
$ws.Range("F10").Value = "rz-user>.header__button.ng-star-inserted"
$ws.Range("F11").Value = "rz-cart>.header__button.ng-star-inserted"
$ws.Range("F12").Value = ".premium-wrapper.ng-star-inserted"
$ws.Range("F13").Value = ".main-slider__pagination-link"
$ws.Range("F14").Value = ".button.button--navy"
$ws.Range("F15").Value = "button[title = 'MasterCard Secure']"
$ws.Range("F16").Value = 'a[title = "Приложение для Андроида"]'
$ws.Range("F17").Value = 'a[title = "Приложение для Айфона"]'
$ws.Range("F4").Value  = 'a[href="https://rozetka.com.ua/cabinet/orders/"]'

# Nudge the screenshot picture a little shorter (its bottom edge moved up
# from row 16 to row 15 in the original edit). Resize before touching row
# heights below, so the anchor is computed against the original layout.
$pic = $ws.Shapes.Item(1)
$pic.Height = 480.51023622047245

# Row 7 now wraps onto four lines inside its styled (wrap-text) cell -
# grow the row to fit, like Excel would when the text is entered.
$ws.Rows.Item(7).RowHeight = 68

# Turn the Precondition URL in B18 into a real hyperlink (adds the
# Hyperlink cell style/font and the <hyperlinks> part).
$cell = $ws.Range("B18")
$ws.Hyperlinks.Add($cell, $cell.Value2) | Out-Null

# Restore the selection to where the author left off.
$ws.Range("G7").Select() | Out-Null
